# Edit: Add temperature distribution data table (rows 12-36) to Sheet1,
# matching the new content introduced in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new text cells in the exact order needed so the shared-string
#     table gets built with the same index sequence as the target workbook ---
$ws.Range("D13").Value = "4 inlets"
$ws.Range("F13").Value = "5 inlets "
$ws.Range("H13").Value = "6 inlets"
$ws.Range("D12").Value = "Temperature distribution at outlet"
$ws.Range("C36").Value = "Std deviation"
$ws.Range("C35").Value = "Mean"
$ws.Range("J13").Value = "No inlets"
$ws.Range("E13").Value = "4 inlets"
$ws.Range("G13").Value = "5 inlets "

# --- Row 12: merge + center the header cell ---
$ws.Range("D12:H12").HorizontalAlignment = -4108
$ws.Range("D12:H12").Merge()

# --- Rows 14-33: temperature distribution data (D:H), 20 rows ---
$data = New-Object 'object[,]' 20,5
$data[0,0] = 566.74871015357701
$data[0,1] = 565.69928417120298
$data[0,2] = 583.76928895483695
$data[0,3] = 574.65451236076399
$data[0,4] = 613.74589739406599
$data[1,0] = 614.01669315346805
$data[1,1] = 610.96006462217599
$data[1,2] = 653.253930527428
$data[1,3] = 634.71303532588297
$data[1,4] = 694.40081850844297
$data[2,0] = 661.36499984816896
$data[2,1] = 656.55504493636204
$data[2,2] = 711.509211630398
$data[2,3] = 690.13229321584197
$data[2,4] = 750.46063614382399
$data[3,0] = 708.86503327397702
$data[3,1] = 702.65770306322497
$data[3,2] = 762.48277435623402
$data[3,3] = 741.58856067334398
$data[3,4] = 796.172770085602
$data[4,0] = 756.62276654734103
$data[4,1] = 749.406911459561
$data[4,2] = 809.46300238863205
$data[4,3] = 790.40096828694902
$data[4,4] = 837.8002550274
$data[5,0] = 795.53018832548901
$data[5,1] = 787.60646442235895
$data[5,2] = 845.901029369924
$data[5,3] = 828.88028819041597
$data[5,4] = 869.80453007930896
$data[6,0] = 826.07167677721498
$data[6,1] = 817.68763685832198
$data[6,2] = 873.52301702804198
$data[6,3] = 858.32702029076802
$data[6,4] = 894.01132527059599
$data[7,0] = 848.63556795512295
$data[7,1] = 839.97615631660699
$data[7,2] = 893.44795726604502
$data[7,3] = 879.68276205147095
$data[7,4] = 911.47784588247896
$data[8,0] = 863.51467739920997
$data[8,1] = 854.70750967469496
$data[8,2] = 906.38602353881004
$data[8,3] = 893.590735744585
$data[8,4] = 922.83072163081795
$data[9,0] = 870.90508243458498
$data[9,1] = 862.03524494547003
$data[9,2] = 912.75731078317597
$data[9,3] = 900.44955396428998
$data[9,4] = 928.42603929275595
$data[10,0] = 870.90507481051702
$data[10,1] = 862.03523792588499
$data[10,2] = 912.75737069802096
$data[10,3] = 900.44959057592905
$data[10,4] = 928.42609821468
$data[11,0] = 863.51465466082698
$data[11,1] = 854.70748909581698
$data[11,2] = 906.386201293873
$data[11,3] = 893.5908434163
$data[11,4] = 922.83089753296997
$data[12,0] = 848.63553061819596
$data[12,1] = 839.97612363479902
$data[12,2] = 893.448246583931
$data[12,3] = 879.68293415626499
$data[12,4] = 911.47813592352202
$data[13,0] = 826.07162607087196
$data[13,1] = 817.68759454102405
$data[13,2] = 873.52340638069199
$data[13,3] = 858.32724537226704
$data[13,4] = 894.01172379603599
$data[14,0] = 795.53012669467398
$data[14,1] = 787.60641595507502
$data[14,2] = 845.90150028197797
$data[14,3] = 828.88054954170002
$data[14,4] = 869.80502693292306
$data[15,0] = 756.62269821539098
$data[15,1] = 749.40686124955096
$data[15,2] = 809.46352695990402
$data[15,3] = 790.40124385014599
$data[15,4] = 837.80083253836597
$data[16,0] = 708.86496463919605
$data[16,1] = 702.65765614097097
$data[16,2] = 762.48331216421195
$data[16,3] = 741.58882387974802
$data[16,4] = 796.17339807574399
$data[17,0] = 661.36494003955897
$data[17,1] = 656.55500685188599
$data[17,2] = 711.50970180790603
$data[17,3] = 690.13251290368498
$data[17,4] = 750.46125887711901
$data[18,0] = 614.01665193553094
$data[18,1] = 610.96003997145101
$data[18,2] = 653.25429488525594
$data[18,3] = 634.71318246173701
$data[18,4] = 694.40134615075601
$data[19,0] = 566.74869536635697
$data[19,1] = 565.69927567898696
$data[19,2] = 583.76943632320695
$data[19,3] = 574.65456513549202
$data[19,4] = 613.74617598095006
$ws.Range("D14:H33").Value = $data

# --- Row 35: Mean ---
$ws.Range("D35").Formula = "=AVERAGE(D14:D33)"
$ws.Range("E35:G35").Formula = "=AVERAGE(E14:E33)"

# --- Row 36: Std deviation ---
$ws.Range("D36").Formula = "=STDEV.S(D14:D33)"
$ws.Range("E36:G36").Formula = "=STDEV.S(E14:E33)"
$ws.Range("H36").Formula = "=STDEV.S(H14:H33)"

# --- Column width adjustment for column C ---
$ws.Range("C1").ColumnWidth = 11.21875

# --- View state: scroll position & selection ---
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("D36:H36").Select()
